# Jesse's Week 6 logs: copy over local activity/task data.
$wb = $excel.ActiveWorkbook

# --- TASK SUMMARY SHEET ---
$ts = $wb.Worksheets.Item("TASK SUMMARY SHEET")

# Header: name + week number
$ts.Range("C1").Value = "Jesse Hare"
$ts.Range("E1").Value = 6

# Task rows
$ts.Range("A3").Value = "Project Build"
$ts.Range("B3").Value = "Begin work on dynamic search/filter"
$ts.Range("C3").Value = 10
$ts.Range("D3").Value = 10
$ts.Range("E3").Value = 6

$ts.Range("A4").Value = "Project Build"
$ts.Range("B4").Value = "Try to reduce memory footprint of program"
$ts.Range("C4").Value = 3
$ts.Range("D4").Value = 4
$ts.Range("E4").Value = 0

$ts.Range("A5").Value = "Project Build"
$ts.Range("B5").Value = "performance analysis and optimisation of program"
$ts.Range("C5").Value = 4
$ts.Range("D5").Value = 3
$ts.Range("E5").Value = 0

$ts.Range("A6").Value = "Project Build"
$ts.Range("B6").Value = "Rewrite to only have one db connection that persists during runtime"
$ts.Range("C6").Value = 1
$ts.Range("D6").Value = 3
$ts.Range("E6").Value = 0

# --- ACTIVITY LOG SUMMARY SHEET ---
$al = $wb.Worksheets.Item("ACTIVITY LOG SUMMARY SHEET")

$al.Range("D1").Value = "Jesse Hare"

$al.Range("A4").Value = "Project Build"
$al.Range("B4").Value = 10
$al.Range("C4").Value = 10
